# Checklist.xlsx - mark LeetCode #152 "Maximum Product Subarray" complete.
#
# The "List" sheet tracks solved problems as rows. A new row is inserted
# right above the #215 entry (which currently sits at row 10) with the
# data for problem #152, and the old rows 10-12 shift down to 11-13.
# The yellow "next up" highlight style that used to live on column A is
# also cleared everywhere, since the tracker has moved past it.

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List")

# Insert a new row above row 10 (the #215 row), shifting it and the rows
# below it down by one.
$ws.Rows.Item(10).Insert()

# Copy the date/note-column formatting from the row above (#148, still at
# row 9) onto the new row before filling in values, so the new row gets
# the same styles (date format, wrapped/centered note cells) instead of
# Excel inventing brand-new style entries.
$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial($xlPasteFormats)
$ws.Range("F9:J9").Copy()
$ws.Range("F10:J10").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# Fill in the new row for problem #152.
$ws.Cells.Item(10, 1).Value = 152
$ws.Cells.Item(10, 2).Value = "Maximum Product Subarray"
$ws.Cells.Item(10, 3).Value = 45091
$ws.Cells.Item(10, 4).Value = "Medium"
$ws.Cells.Item(10, 5).Value = "C"

# Row 10 keeps the default (automatic) row height - no explicit height set.

# Remove the yellow "current" highlight style from column A across all
# data rows; it is no longer used anywhere in the sheet.
for ($r = 3; $r -le 13; $r++) {
  $ws.Cells.Item($r, 1).Style = "Normal"
}

# Restore the usual view: scrolled up a bit with E10 selected.
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("E10").Select() | Out-Null
